# Update crypto price/volume figures per the latest scrape (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.151.68"
$ws.Range("E2").Value = "  -2.19%  "
$ws.Range("D3").Value = "1.851.72"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'237.23"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("E6").Value = "  -5.72%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.07741"
$ws.Range("E8").Value = "  +7.87%  "
$ws.Range("E9").Value = "  -3.14%  "
$ws.Range("D10").Value = "'23.16"
$ws.Range("E10").Value = "  -5.65%  "
$ws.Range("D11").Value = "'0.08157"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "'0.7226"
$ws.Range("E12").Value = "  -3.09%  "
$ws.Range("D13").Value = "1.827.61"
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").Value = "'5.193"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").Value = "'89.34"
$ws.Range("E15").Value = "  -3.65%  "
$ws.Range("D16").Value = "29.149.55"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").Value = "'0.000007811"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'5.727"
$ws.Range("E18").Value = "  -4.92%  "
$ws.Range("D19").Value = "'13.15"
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").Value = "'233.91"
$ws.Range("E20").Value = "  -5.64%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "2.100.99"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "'0.9999"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'7.448"
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").Value = "'161.63"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("D26").Value = "'8.948"
$ws.Range("E26").Value = "  -3.25%  "
$ws.Range("D27").Value = "'0.1430"
$ws.Range("E27").Value = "  -6.45%  "
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "'1.963"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "'1.402"
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("D32").Value = "'1.485"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").Value = "'4.007"
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("D34").Value = "'0.05183"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("D35").Value = "'1.178"
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("D36").Value = "'1.026"
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("D37").Value = "'0.7032"
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("D38").Value = "'2.658"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("D40").Value = "'2.675"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("D41").Value = "'0.9183"
$ws.Range("D42").Value = "1.102.31"
$ws.Range("E42").Value = "  +6.36%  "
$ws.Range("D43").Value = "'6.001"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "'0.4269"
$ws.Range("E44").Value = "  -4.60%  "
$ws.Range("D45").Value = "'69.97"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'102.55"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "'1.756"
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("D49").Value = "1.996.42"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").Value = "'9.168"
$ws.Range("E50").Value = "  -5.04%  "
$ws.Range("D51").Value = "'6.899"
$ws.Range("E51").Value = "  -7.84%  "
